# Applies the "7/13/20 news" update:
#  - adds a new Title paragraph ("EMR News") at the top of the document
#  - rewrites the first (FirstParagraph-styled) news item with the Jan 03, 2020 entry
#  - replaces the trailing horizontal-rule paragraph with a sequence of
#    BodyText-styled news items (Feb 04 .. Jun 30, 2020)

function New-NewsParagraph {
    param(
        [int]$BeforeIndex,
        [string]$StyleName,
        $Runs
    )

    $doc = $word.ActiveDocument

    $anchor = $doc.Paragraphs.Item($BeforeIndex)
    $anchor.Range.InsertParagraphBefore()
    $newPara = $doc.Paragraphs.Item($BeforeIndex)

    # Compose the paragraph's plain text up front (using a vertical-tab
    # placeholder, chr(11), for manual line breaks), remembering which
    # character spans need bold applied afterwards.
    $fullText = ""
    $boldSpans = @()
    foreach ($run in $Runs) {
        $spanStart = $fullText.Length
        if ($run.t -eq "br") {
            $fullText = $fullText + [char]11
        } else {
            $fullText = $fullText + $run.s
        }
        $spanEnd = $fullText.Length
        if ($run.t -eq "text" -and $run.bold) {
            $boldSpans += @{ start = $spanStart; end = $spanEnd }
        }
    }

    # Assign the paragraph text before the style, so that any bold baked
    # into the style definition itself (e.g. the Title style) is not
    # stamped as explicit run formatting.
    $newPara.Range.Text = $fullText
    $newPara.Style = $StyleName
    $paraStart = $newPara.Range.Start

    foreach ($span in $boldSpans) {
        $boldRange = $doc.Range($paraStart + $span.start, $paraStart + $span.end)
        $boldRange.Font.Bold = 1
    }
}

$d = $word.ActiveDocument

# --- New Title paragraph: "EMR News" ---
New-NewsParagraph 1 "Title" @(
        @{ t = "text"; bold = $false; s = ("EMR") },
        @{ t = "text"; bold = $false; s = (" ") },
        @{ t = "text"; bold = $false; s = ("News") }
    )

# --- Rewrite the FirstParagraph news item (was the "NA" / "NA" placeholder) ---
$firstNewsPara = $d.Paragraphs.Item(2)
$dateFind = $firstNewsPara.Range.Find
$dateFind.ClearFormatting()
$dateFind.Execute("NA", $false, $false, $false, $false, $false, $true, 1, $false, "Jan 03, 2020", 1) | Out-Null

$bodyFind = $firstNewsPara.Range.Find
$bodyFind.ClearFormatting()
$bodyFind.Execute("NA", $false, $false, $false, $false, $false, $true, 1, $false, "a global automation technology and engineering company, has been named the" + " " + [char]8220 + "Industrial IoT Company of the Year" + [char]8221 + " " + "by IoT Breakthrough for an unprecedented third consecutive year. The honor recognizes Emerson" + [char]8217 + "s commitment to helping customers in industries such as chemical, life sciences, power, and oil and gas define and execute a practical and successful path to digital transformation. Emerson recently introduced a new, dedicated digital transformation business that combines Emerson" + [char]8217 + "s leading sensing technology, operational analytics and broad services capabilities to deliver targeted digital solutions to customer challenges.", 1) | Out-Null

# --- Replace the horizontal-rule paragraph with the dated BodyText news items ---
$hrIndex = 3
New-NewsParagraph $hrIndex "BodyText" @(
        @{ t = "text"; bold = $true; s = ("Feb 04, 2020") },
        @{ t = "br" },
        @{ t = "text"; bold = $false; s = ("First-quarter earnings were `$326 million, or 53 cents a share, down 30% compared with `$465 million, or 74 cents a share a year earlier. Adjusted earnings were 67 cents, matching the expectations of analysts polled by FactSet. Total revenue was `$4.15 billion, remaining roughly the same from a year earlier. Analysts were expecting revenue of `$4.15 billion. Emerson expects FY2020 adjusted EPS of `$3.55 to `$3.80 compared with its previous outlook of `$3.48 to `$3.72 and the Street projection of `$3.63. The company said the improved outlook reflects the favorable impact of expected savings from its restructuring efforts.") }
    )
$hrIndex = $hrIndex + 1
New-NewsParagraph $hrIndex "BodyText" @(
        @{ t = "text"; bold = $true; s = ("Feb 14, 2020") },
        @{ t = "br" },
        @{ t = "text"; bold = $false; s = ("Says it has no plans to break up its current portfolio of businesses. During a presentation to investors the industrial conglomerate says there will be") },
        @{ t = "text"; bold = $false; s = (" ") },
        @{ t = "text"; bold = $false; s = ([char]8220) },
        @{ t = "text"; bold = $false; s = ("no breakup unless a major strategic acquisition") },
        @{ t = "text"; bold = $false; s = ([char]8221) },
        @{ t = "text"; bold = $false; s = (" ") },
        @{ t = "text"; bold = $false; s = ("is undertaken that would require a move.") }
    )
$hrIndex = $hrIndex + 1
New-NewsParagraph $hrIndex "BodyText" @(
        @{ t = "text"; bold = $true; s = ("Feb 18, 2020") },
        @{ t = "br" },
        @{ t = "text"; bold = $false; s = ("RBC Capital Markets analyst Deane Dray cut the recommendation on Emerson Electric Co. to sector perform from outperform. PT set to `$74, implies a 0.7% increase from last price.") }
    )
$hrIndex = $hrIndex + 1
New-NewsParagraph $hrIndex "BodyText" @(
        @{ t = "text"; bold = $true; s = ("Feb 21, 2020") },
        @{ t = "br" },
        @{ t = "text"; bold = $false; s = ("Was cut to sector perform from outperform at RBC Capital Markets, which wrote that its") },
        @{ t = "text"; bold = $false; s = (" ") },
        @{ t = "text"; bold = $false; s = ([char]8220) },
        @{ t = "text"; bold = $false; s = ("bull case for a breakup appears to be postponed.") },
        @{ t = "text"; bold = $false; s = ([char]8221) }
    )
$hrIndex = $hrIndex + 1
New-NewsParagraph $hrIndex "BodyText" @(
        @{ t = "text"; bold = $true; s = ("Feb 28, 2020") },
        @{ t = "br" },
        @{ t = "text"; bold = $false; s = ("Estimates 2Q sales impact from coronavirus at least `$100m to `$150m. Company had estimated at investor conference on Feb. 13 that impact would be `$75m to `$100m. Emerson now says half of the expected 2Q sales impact expected to be recovered in the fiscal year, but cautions that the longer the virus endures, the lower the likelihood of sales recovery this year.") }
    )
$hrIndex = $hrIndex + 1
New-NewsParagraph $hrIndex "BodyText" @(
        @{ t = "text"; bold = $true; s = ("Mar 13, 2020") },
        @{ t = "br" },
        @{ t = "text"; bold = $false; s = ("Wolfe Research analyst Nigel Coe cut the recommendation on Emerson Electric Co. to peerperform from outperform.Price target set to `$71, implies a 48% increase from last price. Emerson Electric average price target is `$77.60.") }
    )
$hrIndex = $hrIndex + 1
New-NewsParagraph $hrIndex "BodyText" @(
        @{ t = "text"; bold = $true; s = ("Apr 21, 2020") },
        @{ t = "br" },
        @{ t = "text"; bold = $false; s = ("The company reported sales fell 9% to `$4.16 billion in its fiscal second quarter, down from `$4.57 billion a year earlier. Analysts polled by FactSet had forecast `$4.28 billion in sales for the latest period. Profit slipped to `$517 million, or 84 cents a share, from `$520 million, or 84 cents a share, the year earlier. Emerson" + [char]8217 + "s adjusted profit of 89 cents a share beat the consensus estimate for that metric by 12 cents. Lower costs helped results. Selling, general and administrative expenses dropped to `$983 million from `$1.15 billion. Moving forward, Emerson is anticipating FY2020 adjusted EPS of `$3 to `$3.20 compared with the previous projected range of `$3.55 to `$3.80. The company said that 2020 buybacks will be approximately `$950 million, compared with its prior guidance of `$1.5 billion, and capital expenditures are projected to be `$550 million, down from the previous estimate of `$650 million. It also maintained its current dividend policy.") }
    )
$hrIndex = $hrIndex + 1
New-NewsParagraph $hrIndex "BodyText" @(
        @{ t = "text"; bold = $true; s = ("Jun 30, 2020") },
        @{ t = "br" },
        @{ t = "text"; bold = $false; s = ("Price Target Raised to `$59.00/Share From `$54.00 by Morgan Stanley.") }
    )
$hrIndex = $hrIndex + 1

# Remove the now-redundant horizontal-rule paragraph.
$hrPara = $d.Paragraphs.Item($hrIndex)
$hrPara.Range.Delete()

Write-Host "Final paragraph count:" $d.Paragraphs.Count
